$wb = $excel.ActiveWorkbook

# --- Sheet2: remove the "12/4 362p" position row (row 58) ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Rows("58:58").Delete()

# --- Sheet1: update the manually-tallied calls/puts counts that changed
#     because that position was removed ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C88").Value = 31
$ws1.Range("C89").Value = 73
$ws1.Range("C99").Value = 1
$ws1.Range("C135").Value = 4
$ws1.Range("C139").Value = 8
$ws1.Range("C177").Value = 31
$ws1.Range("C178").Value = 39
